$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.402.42"
$ws.Range("E2").Value = "'  +0.02%  "
$ws.Range("D3").Value = "'1.718.00"
$ws.Range("E3").Value = "'  +0.11%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  +1.06%  "
$ws.Range("D5").Value = "'224.64"
$ws.Range("E5").Value = "'  -3.10%  "
$ws.Range("D6").Value = "'0.5322"
$ws.Range("E6").Value = "'  -2.40%  "
$ws.Range("E7").Value = "'  +1.03%  "
$ws.Range("D8").Value = "'0.2645"
$ws.Range("E8").Value = "'  -4.03%  "
$ws.Range("D9").Value = "'0.06555"
$ws.Range("E9").Value = "'  +0.57%  "
$ws.Range("D10").Value = "'21.17"
$ws.Range("E10").Value = "'  -2.99%  "
$ws.Range("D11").Value = "'0.07658"
$ws.Range("E11").Value = "'  -0.98%  "
$ws.Range("D12").Value = "'4.585"
$ws.Range("E12").Value = "'  -3.56%  "
$ws.Range("B13").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "'1.957.59"
$ws.Range("E13").Value = "'  +0.50%  "
$ws.Range("B14").Value = "'WrappedEther"
$ws.Range("C14").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.696.93"
$ws.Range("E14").Value = "'  -0.98%  "
$ws.Range("D15").Value = "'0.5752"
$ws.Range("E15").Value = "'  -5.80%  "
$ws.Range("D16").Value = "'0.0₅8211"
$ws.Range("E16").Value = "'  -2.01%  "
$ws.Range("D17").Value = "'67.51"
$ws.Range("E17").Value = "'  -2.35%  "
$ws.Range("D18").Value = "'27.414.86"
$ws.Range("E18").Value = "'  +0.36%  "
$ws.Range("D19").Value = "'216.77"
$ws.Range("E19").Value = "'  +2.11%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "'  +0.99%  "
$ws.Range("D21").Value = "'4.705"
$ws.Range("E21").Value = "'  -2.21%  "
$ws.Range("D22").Value = "'10.51"
$ws.Range("E22").Value = "'  -5.08%  "
$ws.Range("D23").Value = "'5.960"
$ws.Range("E23").Value = "'  -4.55%  "
$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "'  +1.34%  "
$ws.Range("D25").Value = "'143.00"
$ws.Range("E25").Value = "'  -2.88%  "
$ws.Range("D26").Value = "'1.736"
$ws.Range("E26").Value = "'  +7.78%  "
$ws.Range("D27").Value = "'0.1222"
$ws.Range("E27").Value = "'  -3.42%  "
$ws.Range("D28").Value = "'7.282"
$ws.Range("E28").Value = "'  -2.35%  "
$ws.Range("D29").Value = "'16.38"
$ws.Range("E29").Value = "'  -2.36%  "
$ws.Range("D30").Value = "'0.05412"
$ws.Range("E30").Value = "'  -4.92%  "
$ws.Range("D31").Value = "'1.295"
$ws.Range("E31").Value = "'  -2.28%  "
$ws.Range("D32").Value = "'3.507"
$ws.Range("E32").Value = "'  -4.29%  "
$ws.Range("D33").Value = "'3.422"
$ws.Range("E33").Value = "'  -2.89%  "
$ws.Range("D34").Value = "'1.639"
$ws.Range("E34").Value = "'  -0.42%  "
$ws.Range("D35").Value = "'2.876"
$ws.Range("E35").Value = "'  +0.37%  "
$ws.Range("B36").Value = "'HuobiToken"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.434"
$ws.Range("E36").Value = "'  +0.70%  "
$ws.Range("B37").Value = "'ARBITRUM"
$ws.Range("C37").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "'0.9505"
$ws.Range("E37").Value = "'  -3.46%  "
$ws.Range("D38").Value = "'0.5876"
$ws.Range("E38").Value = "'  -0.13%  "
$ws.Range("D39").Value = "'0.01636"
$ws.Range("E39").Value = "'  -0.99%  "
$ws.Range("D40").Value = "'5.883"
$ws.Range("E40").Value = "'  -0.87%  "
$ws.Range("B41").Value = "'PaxDollar"
$ws.Range("C41").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "'1.006"
$ws.Range("E41").Value = "'  +1.09%  "
$ws.Range("B42").Value = "'Maker"
$ws.Range("C42").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "'1.042.37"
$ws.Range("E42").Value = "'  -1.59%  "
$ws.Range("B43").Value = "'TrustWalletToken"
$ws.Range("C43").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8424"
$ws.Range("E43").Value = "'  +0.23%  "
$ws.Range("D44").Value = "'101.15"
$ws.Range("E44").Value = "'  -1.80%  "
$ws.Range("D45").Value = "'1.863.80"
$ws.Range("E45").Value = "'  +0.45%  "
$ws.Range("D46").Value = "'0.0₈115"
$ws.Range("E46").Value = "'  +8.19%  "
$ws.Range("D47").Value = "'58.26"
$ws.Range("E47").Value = "'  -2.98%  "
$ws.Range("D48").Value = "'0.4516"
$ws.Range("E48").Value = "'  +4.17%  "
$ws.Range("B49").Value = "'XinFinNetwork"
$ws.Range("C49").Value = "'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D49").Value = "'0.06658"
$ws.Range("E49").Value = "'  +16.00%  "
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.155"
$ws.Range("E50").Value = "'  +0.05%  "
$ws.Range("B51").Value = "'Frax"
$ws.Range("C51").Value = "'https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "'  +0.25%  "
